# "update Time $ Accuracy"
# Refresh the algorithm comparison table: relabel the algorithms, refresh the
# timing/accuracy figures, reorder "Deep Learning" to the bottom of the list,
# and give the table a bordered/shaded look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ------------------------------------------------------
$ws.Range("A1").Value = "Algorithm"
$ws.Range("B1").Value = "Time"
$ws.Range("C1").Value = "Accuracy"

# ---- Data rows (Deep Learning moved from the top to the bottom) ------
$ws.Range("A2").Value = "Random Forest"
$ws.Range("B2").Value = "6.558833s"
$ws.Range("C2").Value = 0.91395

$ws.Range("A3").Value = "SVM"
$ws.Range("B3").Value = "5.836301s"
$ws.Range("C3").Value = 0.9167

$ws.Range("A4").Value = "AdaBoost"
$ws.Range("B4").Value = "5.571954s"
$ws.Range("C4").Value = 0.9352

$ws.Range("A5").Value = "Deep Learning"
$ws.Range("B5").Value = "5.430689s"
$ws.Range("C5").Value = 0.96

# ---- Formatting: bordered cells with a white fill, black text --------
$table = $ws.Range("A1:C5")
$table.Interior.Color = 16777215
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2
$table.HorizontalAlignment = -4108
$table.VerticalAlignment = -4108

$header = $ws.Range("A1:C1")
$header.Font.Color = 0
$header.Font.Bold = $true

$data = $ws.Range("A2:C5")
$data.Font.Color = 0
$data.Font.Bold = $false

# ---- Selection, matching the saved workbook state ---------------------
$ws.Range("B7").Select() | Out-Null
